# Update countries & provincias Spain
# Refresh country COVID stats table (sheet "Pais") and the "last updated" timestamp,
# and fix the Sri Lanka / Cabo Verde row ordering/labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Agosto de 2020 a las 21:10"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 5173322
$ws.Range("C4").Value = 23599
$ws.Range("D4").Value = 2640740
$ws.Range("E4").Value = 2367285
$ws.Range("G4").Value = 227
$ws.Range("H4").Value = 165297

# --- Row 6: India ---
$ws.Range("B6").Value = 2214137
$ws.Range("C6").Value = 62117
$ws.Range("D6").Value = 1534278
$ws.Range("E6").Value = 635393
$ws.Range("G6").Value = 1013
$ws.Range("H6").Value = 44466

# --- Row 21: Turquia ---
$ws.Range("B21").Value = 240804
$ws.Range("C21").Value = 1182
$ws.Range("D21").Value = 223759
$ws.Range("E21").Value = 11201
$ws.Range("G21").Value = 15
$ws.Range("H21").Value = 5844

# --- Row 31: Ecuador ---
$ws.Range("E31").Value = 16932
$ws.Range("G31").Value = 6
$ws.Range("H31").Value = 5922

# --- Row 62: Suiza ---
$ws.Range("B62").Value = 30609
$ws.Range("C62").Value = 957
$ws.Range("D62").Value = 22042
$ws.Range("E62").Value = 8373
$ws.Range("G62").Value = 7
$ws.Range("H62").Value = 194

# --- Row 79: Estado de Palestina ---
$ws.Range("B79").Value = 14208
$ws.Range("C79").Value = 280
$ws.Range("E79").Value = 6167

# --- Rows 120/121: Sri Lanka <-> Cabo Verde swap with refreshed data ---
$ws.Range("A120").Value = "Cabo Verde"
$ws.Range("B120").Value = 2858
$ws.Range("C120").Value = 23
$ws.Range("D120").Value = 2086
$ws.Range("E120").Value = 740
$ws.Range("H120").Value = 32

$ws.Range("A121").Value = "Sri Lanka"
$ws.Range("B121").Value = 2841
$ws.Range("D121").Value = 2579
$ws.Range("E121").Value = 251
$ws.Range("H121").Value = 11
